# DAS-677 - CCRU - Creation of Scenes for SOVI SOCVI
#
# The original "Panoramic Photo" / "Panoramic photo of Cooler" / "Menu Board,
# Cash Zone, SS_Menu Board, SS_Cash Zone" tagging values (column Y / X on the
# "FF" sheet) are extended with the new Scene-tagging aliases used for the
# SOVI/SOCVI scenes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FF")

# Column Y ("Scenes to include") on the SKU-availability rows: append the
# "SS_Panoramic Photo" alias.
$panoramicPhotoRows = @(4,5,6,7,8,9,10,11,12,13,14,15,16,18,19,21,22,24,25,26,28,29,30,31)
foreach ($r in $panoramicPhotoRows) {
    $ws.Cells.Item($r, 25).Value = "Panoramic Photo, SS_Panoramic Photo"
}

# Column X ("Scenes to exclude") row 38 ("Cooler: Prime Position"): append the
# "SS_Panoramic photo of Cooler" alias.
$ws.Cells.Item(38, 24).Value = "Panoramic photo of Cooler, SS_Panoramic photo of Cooler"

# Column X rows 42-45 (Menu Activation / Combo rows): extend the Menu
# Board / Cash Zone scene list with the new Canteen/QSR scene sub-types.
$menuBoardRows = @(42,43,44,45)
foreach ($r in $menuBoardRows) {
    $ws.Cells.Item($r, 24).Value = "Menu Board, Cash Zone, SS_Cash Zone - Canteen, QSR, SS_Menu Board - Canteen, QSR"
}

# Refresh the window view: unfreeze scroll back to the top of the scrollable
# area (below the frozen header row) and move the active selection to the
# newly-edited X43:X45 block.
$ws.Range("X43:X45").Select()
